# Fruta / hortaliza, semanal
# Insert a new weekly record at row 126 of the "Piña" sheet, pushing the
# existing rows 126-205 down to 127-206 (dimension grows from A1:T205 to
# A1:T206).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 126; Excel shifts rows 126:205
# down to 127:206 and carries the row-above formatting (incl. the date
# number format on column D) onto the newly inserted row.
$ws.Rows("126").Insert()

# Populate the newly inserted row 126 with the new weekly observation.
$ws.Range("A126").Value = 5
$ws.Range("B126").Value = "Macroferia Regional de Talca"
$ws.Range("C126").Value = "Maule"
$ws.Range("D126").Value = 44596
$ws.Range("E126").Value = 7
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100108
$ws.Range("H126").Value = "Tropicales y subtropicales"
$ws.Range("I126").Value = 100108005
$ws.Range("J126").Value = "Piña"
$ws.Range("K126").Value = "Caramelo"
$ws.Range("L126").Value = "Tercera"
$ws.Range("M126").Value = 200
$ws.Range("N126").Value = 16000
$ws.Range("O126").Value = 16000
$ws.Range("P126").Value = 16000
$ws.Range("Q126").Value = "$/caja 16 unidades"
$ws.Range("R126").Value = "Ecuador"
$ws.Range("S126").Value = 1000
$ws.Range("T126").Value = 16
